$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.4165500938170587
$ws.Range("C2").Value = 0.5817745553736098
$ws.Range("D2").Value = 0.6656782362548082
$ws.Range("E2").Value = 0.7263964236600271
$ws.Range("B3").Value = 0.4425011968473371
$ws.Range("C3").Value = 0.6127678567160914
$ws.Range("D3").Value = 0.6945291719047449
$ws.Range("E3").Value = 0.7497963153137237
$ws.Range("B4").Value = 0.3717759532087446
$ws.Range("C4").Value = 0.5267388280248954
$ws.Range("D4").Value = 0.6115471694156789
$ws.Range("E4").Value = 0.6784402355552254
$ws.Range("B5").Value = 0.4802277040075707
$ws.Range("C5").Value = 0.6573702918581965
$ws.Range("D5").Value = 0.7293867604872588
$ws.Range("E5").Value = 0.7845089806605176
$ws.Range("B6").Value = 0.4707275492917793
$ws.Range("C6").Value = 0.6459004604667629
$ws.Range("D6").Value = 0.7187228668629096
$ws.Range("E6").Value = 0.7756754437877884
$ws.Range("B7").Value = 0.4945077918741457
$ws.Range("C7").Value = 0.6728827465569309
$ws.Range("D7").Value = 0.7498245758194799
$ws.Range("E7").Value = 0.7937543922314764
$ws.Range("B8").Value = 0.4303344470155975
$ws.Range("C8").Value = 0.6002100318368492
$ws.Range("D8").Value = 0.6903222019640344
$ws.Range("E8").Value = 0.7384478743081042
$ws.Range("B9").Value = 0.4982263498511468
$ws.Range("C9").Value = 0.6793750769462492
$ws.Range("D9").Value = 0.7568206031162165
$ws.Range("E9").Value = 0.8003621903725986
$ws.Range("B10").Value = 0.5123575806189912
$ws.Range("C10").Value = 0.6874373528801301
$ws.Range("D10").Value = 0.7523855043214686
$ws.Range("E10").Value = 0.7818163400102417
$ws.Range("B11").Value = 0.5075664387102803
$ws.Range("C11").Value = 0.6805110419387205
$ws.Range("D11").Value = 0.7447263700905355
$ws.Range("E11").Value = 0.7742004243341288
$ws.Range("B12").Value = 0.4733464200831813
$ws.Range("C12").Value = 0.6156258578351426
$ws.Range("D12").Value = 0.651220674327408
$ws.Range("E12").Value = 0.6573517222229365
$ws.Range("B13").Value = 0.5100629818656438
$ws.Range("C13").Value = 0.6837280873050913
$ws.Range("D13").Value = 0.7479563703930113
$ws.Range("E13").Value = 0.77706033667528
